# Updated cryptos list on Sun Jul  9 12:48:13 UTC 2023 with GitHub Actions
#
# The "Price" (D) and "Volume(1h)" (E) columns hold values that look
# numeric ("1.000", "0.9996", "  +0.28%  ", ...) but must stay literal
# text, exactly as they were authored (matching the original t="inlineStr"
# cells). Plain `.Value = "1.000"` assignment lets Excel's COM layer
# auto-coerce numeric-looking strings into real numbers (e.g. "1.000" -> 1),
# so we briefly force the destination range to Text format before writing,
# then restore the "Normal" style afterwards so no stray number-format/style
# index is left behind on the cells (matching the unchanged `s` attributes
# in the original workbook).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$priceVolRange = $ws.Range("D2:E51")
$priceVolRange.NumberFormat = "@"

$ws.Range("D2").Value = "30.301.60"
$ws.Range("E2").Value = "  +0.28%  "

$ws.Range("D3").Value = "1.876.03"
$ws.Range("E3").Value = "  +0.91%  "

$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").Value = "235.11"
$ws.Range("E5").Value = "  -0.82%  "

$ws.Range("D6").Value = "0.9996"
$ws.Range("E6").Value = "  -0.10%  "

$ws.Range("D7").Value = "0.4699"
$ws.Range("E7").Value = "  +0.48%  "

$ws.Range("D8").Value = "0.2879"
$ws.Range("E8").Value = "  +0.48%  "

$ws.Range("D9").Value = "0.06606"
$ws.Range("E9").Value = "  +1.11%  "

$ws.Range("D10").Value = "21.72"
$ws.Range("E10").Value = "  -0.93%  "

$ws.Range("D11").Value = "0.07964"
$ws.Range("E11").Value = "  +0.27%  "

$ws.Range("D12").Value = "96.79"
$ws.Range("E12").Value = "  -0.37%  "

$ws.Range("D13").Value = "1.873.91"
$ws.Range("E13").Value = "  +0.69%  "

$ws.Range("D14").Value = "0.6980"
$ws.Range("E14").Value = "  +2.58%  "

$ws.Range("D15").Value = "5.119"
$ws.Range("E15").Value = "  -1.15%  "

$ws.Range("D16").Value = "270.29"
$ws.Range("E16").Value = "  +1.11%  "

$ws.Range("D17").Value = "30.317.45"
$ws.Range("E17").Value = "  +0.38%  "

$ws.Range("E18").Value = "  +3.18%  "

$ws.Range("D19").Value = "0.000007770"
$ws.Range("E19").Value = "  +5.39%  "

$ws.Range("D20").Value = "0.9992"
$ws.Range("E20").Value = "  -0.10%  "

$ws.Range("D21").Value = "2.119.81"
$ws.Range("E21").Value = "  +0.27%  "

$ws.Range("D22").Value = "0.9999"
$ws.Range("E22").Value = "  -0.12%  "

$ws.Range("D23").Value = "5.278"
$ws.Range("E23").Value = "  -0.82%  "

$ws.Range("D24").Value = "6.222"
$ws.Range("E24").Value = "  +0.38%  "

$ws.Range("D25").Value = "9.398"
$ws.Range("E25").Value = "  +1.94%  "

$ws.Range("D26").Value = "167.82"
$ws.Range("E26").Value = "  +0.36%  "

$ws.Range("D27").Value = "18.94"
$ws.Range("E27").Value = "  +0.30%  "

$ws.Range("D28").Value = "1.957"
$ws.Range("E28").Value = "  +0.29%  "

$ws.Range("E29").Value = "  -1.42%  "

$ws.Range("D30").Value = "0.09886"
$ws.Range("E30").Value = "  +0.44%  "

$ws.Range("D31").Value = "4.349"
$ws.Range("E31").Value = "  -0.43%  "

$ws.Range("D32").Value = "1.464"
$ws.Range("E32").Value = "  -1.39%  "

$ws.Range("D33").Value = "4.065"
$ws.Range("E33").Value = "  +0.22%  "

$ws.Range("D34").Value = "0.04740"
$ws.Range("E34").Value = "  +0.45%  "

$ws.Range("D35").Value = "1.138"
$ws.Range("E35").Value = "  +0.68%  "

$ws.Range("D36").Value = "0.7041"
$ws.Range("E36").Value = "  +0.38%  "

$ws.Range("D37").Value = "2.722"
$ws.Range("E37").Value = "  +0.53%  "

$ws.Range("D38").Value = "0.01877"
$ws.Range("E38").Value = "  -0.01%  "

$ws.Range("D39").Value = "2.813"
$ws.Range("E39").Value = "  +6.98%  "

$ws.Range("D40").Value = "6.236"
$ws.Range("E40").Value = "  -0.11%  "

$ws.Range("D41").Value = "72.08"
$ws.Range("E41").Value = "  -4.17%  "

$ws.Range("D42").Value = "1.963"
$ws.Range("E42").Value = "  +0.90%  "

$ws.Range("D43").Value = "0.4185"
$ws.Range("E43").Value = "  +0.53%  "

# Rows 44/45 swap Coin name + Link (TrustWalletToken now ranks above PaxDollar)
# and also get fresh price/volume figures.
$ws.Range("B44").Value = "TrustWalletToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D44").Value = "0.8427"
$ws.Range("E44").Value = "  -0.85%  "

$ws.Range("B45").Value = "PaxDollar"
$ws.Range("C45").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D45").Value = "0.9989"
$ws.Range("E45").Value = "  -0.10%  "

$ws.Range("D46").Value = "102.61"
$ws.Range("E46").Value = "  -0.61%  "

$ws.Range("D47").Value = "7.137"
$ws.Range("E47").Value = "  -0.46%  "

$ws.Range("D48").Value = "9.132"
$ws.Range("E48").Value = "  -0.92%  "

$ws.Range("D49").Value = "919.81"
$ws.Range("E49").Value = "  -3.82%  "

$ws.Range("D50").Value = "34.67"
$ws.Range("E50").Value = "  +1.63%  "

$ws.Range("D51").Value = "0.05688"
$ws.Range("E51").Value = "  +0.67%  "

# Restore default styling (no explicit number format) now that every value
# in the range is committed as text, so we don't leave a stray style index
# on cells that were style-less in the original file.
$priceVolRange.Style = "Normal"
